# Fill in the newly-found PI controller gains (Kp/Ki derived values) for
# the remaining operating points, then leave the selection where the
# author last left it (I20) before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (speed = 10): only B2/C2 were missing
$ws.Range("B2").Value = 25320
$ws.Range("C2").Value = 11396

# Row 3 (speed = 20): B3:E3 were missing
$ws.Range("B3").Value = 25320
$ws.Range("C3").Value = 11396
$ws.Range("D3").Value = 8300
$ws.Range("E3").Value = 12000

# Row 4 (speed = 40): B4:E4 were missing
$ws.Range("B4").Value = 25320
$ws.Range("C4").Value = 11396
$ws.Range("D4").Value = 8300
$ws.Range("E4").Value = 12000

# Row 5 (speed = 60): B5:E5 were missing
$ws.Range("B5").Value = 25320
$ws.Range("C5").Value = 11396
$ws.Range("D5").Value = 8300
$ws.Range("E5").Value = 12000

# Row 6 (speed = 90): B6:E6 were missing
$ws.Range("B6").Value = 25320
$ws.Range("C6").Value = 11396
$ws.Range("D6").Value = 8300
$ws.Range("E6").Value = 12000

# Move the active selection to I20, matching the saved workbook state.
$ws.Range("I20").Select()
